$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the LOT section (row 2, G:I) to the new "pop_filter2_section" names
$ws.Range("G2").Value = "pop_filter2_section1"
$ws.Range("H2").Value = "pop_filter2_section1_checkbox"
$ws.Range("I2").Value = "pop_filter2_section"

# Rename the Subpopulation section (row 7, G:I) to the new "pop_filter1_section" names
$ws.Range("G7").Value = "pop_filter1_section1"
$ws.Range("H7").Value = "pop_filter1_section1_checkbox"
$ws.Range("I7").Value = "pop_filter1_section"

# Update the selection / top-left cell to match the final view state
$ws.Range("G7:I7").Select()
